# The "RO & CO Hearing Allocation" sheet had a dedicated "Central Office"
# allocation row (row 4, directly under the column headers) that no longer
# belongs in the template. Remove that entire row so every row below it
# shifts up by one (data, styles and row heights all move together), and
# the now-unused "Central Office" shared string is dropped automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RO & CO Hearing Allocation")
$ws.Rows.Item(4).Delete()
